# TC08 - Verify total
# Adds a new worksheet "TC08" at the end of the workbook, populates it with
# the "verify total" test-case data, updates the TC06 selection, and makes
# TC08 the active sheet/tab (matching the commit's xlsx diff).

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet and move it to the end of the tab strip -------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "TC08"
$newSheet.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# --- TC06: move the cell selection to E2 (no change to its data) ---------
$ws06 = $wb.Worksheets.Item("TC06")
$ws06.Range("E2").Select()

# --- TC08: fill in header row + data row ----------------------------------
$ws08 = $wb.Worksheets.Item("TC08")

$ws08.Cells.Item(1, 1).Value = "baseUrl"
$ws08.Cells.Item(1, 2).Value = "userEmail"
$ws08.Cells.Item(1, 3).Value = "password"
$ws08.Cells.Item(1, 4).Value = "size"
$ws08.Cells.Item(1, 5).Value = "quantity"
$ws08.Cells.Item(1, 6).Value = "incressQuantity"

$ws08.Cells.Item(2, 1).Value = "http://automationpractice.com/index.php"
$ws08.Cells.Item(2, 2).Value = "test0000@test.io"
$ws08.Cells.Item(2, 3).Value = "Pass1234"
$ws08.Cells.Item(2, 4).Value = "M"
$ws08.Cells.Item(2, 5).Value = 1
$ws08.Cells.Item(2, 6).Value = 2

# Select F2 and make TC08 the active sheet/tab (activeTab = 7, 0-based) ---
$ws08.Range("F2").Select()
$ws08.Activate()
